$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "F1057-GAGTTGTACG"
$ws.Cells.Item(2, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGAGTTGTACGTCGTCGGCAGCGTC"
$ws.Cells.Item(3, 2).Value = "F1058-CCAGTACGTA"
$ws.Cells.Item(3, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCCAGTACGTATCGTCGGCAGCGTC"
$ws.Cells.Item(4, 2).Value = "F1059-GACACATGAA"
$ws.Cells.Item(4, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGACACATGAATCGTCGGCAGCGTC"
$ws.Cells.Item(5, 2).Value = "F1060-GTCGTAGATG"
$ws.Cells.Item(5, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTCGTAGATGTCGTCGGCAGCGTC"
$ws.Cells.Item(6, 2).Value = "F1061-CAACTCTGTA"
$ws.Cells.Item(6, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCAACTCTGTATCGTCGGCAGCGTC"
$ws.Cells.Item(7, 2).Value = "F1062-AGACTTCCTT"
$ws.Cells.Item(7, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAGACTTCCTTTCGTCGGCAGCGTC"
$ws.Cells.Item(8, 2).Value = "F1063-ATCCACCAAG"
$ws.Cells.Item(8, 3).Value = "AATGATACGGCGACCACCGAGATCTACACATCCACCAAGTCGTCGGCAGCGTC"
$ws.Cells.Item(9, 2).Value = "F1064-TAGTCTGGAA"
$ws.Cells.Item(9, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTAGTCTGGAATCGTCGGCAGCGTC"
$ws.Cells.Item(10, 2).Value = "F1065-TCCTCTACGT"
$ws.Cells.Item(10, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCCTCTACGTTCGTCGGCAGCGTC"
$ws.Cells.Item(11, 2).Value = "F1066-ATGCTCTTGA"
$ws.Cells.Item(11, 3).Value = "AATGATACGGCGACCACCGAGATCTACACATGCTCTTGATCGTCGGCAGCGTC"
$ws.Cells.Item(12, 2).Value = "F1067-TGCTAGACTA"
$ws.Cells.Item(12, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGCTAGACTATCGTCGGCAGCGTC"
$ws.Cells.Item(13, 2).Value = "F1068-TACCTACAGC"
$ws.Cells.Item(13, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTACCTACAGCTCGTCGGCAGCGTC"
$ws.Cells.Item(14, 2).Value = "F1069-AGAAGCAGAG"
$ws.Cells.Item(14, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAGAAGCAGAGTCGTCGGCAGCGTC"
$ws.Cells.Item(15, 2).Value = "F1070-AAGCTTGCAT"
$ws.Cells.Item(15, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAAGCTTGCATTCGTCGGCAGCGTC"
$ws.Cells.Item(16, 2).Value = "F1071-GATCAACATC"
$ws.Cells.Item(16, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGATCAACATCTCGTCGGCAGCGTC"
$ws.Cells.Item(17, 2).Value = "F1072-TGAAGGTTGG"
$ws.Cells.Item(17, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGAAGGTTGGTCGTCGGCAGCGTC"
$ws.Cells.Item(18, 2).Value = "F1073-AGGACAAGGT"
$ws.Cells.Item(18, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAGGACAAGGTTCGTCGGCAGCGTC"
$ws.Cells.Item(19, 2).Value = "F1074-ACGAAGTCTC"
$ws.Cells.Item(19, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACGAAGTCTCTCGTCGGCAGCGTC"
$ws.Cells.Item(20, 2).Value = "F1075-GCAAGACTGT"
$ws.Cells.Item(20, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGCAAGACTGTTCGTCGGCAGCGTC"
$ws.Cells.Item(21, 2).Value = "F1076-ACCATGTCCT"
$ws.Cells.Item(21, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACCATGTCCTTCGTCGGCAGCGTC"
$ws.Cells.Item(22, 2).Value = "F1077-TACTGTGAAG"
$ws.Cells.Item(22, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTACTGTGAAGTCGTCGGCAGCGTC"
$ws.Cells.Item(23, 2).Value = "F1078-TTCGAGTTCC"
$ws.Cells.Item(23, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTTCGAGTTCCTCGTCGGCAGCGTC"
$ws.Cells.Item(24, 2).Value = "F1079-TACATGCTTG"
$ws.Cells.Item(24, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTACATGCTTGTCGTCGGCAGCGTC"
$ws.Cells.Item(25, 2).Value = "F1080-CTACACAACA"
$ws.Cells.Item(25, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTACACAACATCGTCGGCAGCGTC"
$ws.Cells.Item(26, 2).Value = "F1081-GAGAACGTTG"
$ws.Cells.Item(26, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGAGAACGTTGTCGTCGGCAGCGTC"
$ws.Cells.Item(27, 2).Value = "F1082-ACTACACGTA"
$ws.Cells.Item(27, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACTACACGTATCGTCGGCAGCGTC"
$ws.Cells.Item(28, 2).Value = "F1083-ACTGATCGTG"
$ws.Cells.Item(28, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACTGATCGTGTCGTCGGCAGCGTC"
$ws.Cells.Item(29, 2).Value = "F1084-GTACTCCAGT"
$ws.Cells.Item(29, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTACTCCAGTTCGTCGGCAGCGTC"
$ws.Cells.Item(30, 2).Value = "F1085-GTCTCAGTAG"
$ws.Cells.Item(30, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTCTCAGTAGTCGTCGGCAGCGTC"
$ws.Cells.Item(31, 2).Value = "F1086-TCACCTCATC"
$ws.Cells.Item(31, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCACCTCATCTCGTCGGCAGCGTC"
$ws.Cells.Item(32, 2).Value = "F1087-GTTGCAGAGT"
$ws.Cells.Item(32, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTTGCAGAGTTCGTCGGCAGCGTC"
$ws.Cells.Item(33, 2).Value = "F1088-CAACATCCAG"
$ws.Cells.Item(33, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCAACATCCAGTCGTCGGCAGCGTC"
$ws.Cells.Item(34, 2).Value = "F1089-CTACAGATGC"
$ws.Cells.Item(34, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTACAGATGCTCGTCGGCAGCGTC"
$ws.Cells.Item(35, 2).Value = "F1090-GTGCTACGTA"
$ws.Cells.Item(35, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTGCTACGTATCGTCGGCAGCGTC"
$ws.Cells.Item(36, 2).Value = "F1091-TTGTAGCTTC"
$ws.Cells.Item(36, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTTGTAGCTTCTCGTCGGCAGCGTC"
$ws.Cells.Item(37, 2).Value = "F1092-TTCGACAGTC"
$ws.Cells.Item(37, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTTCGACAGTCTCGTCGGCAGCGTC"
$ws.Cells.Item(38, 2).Value = "F1093-CAGCAAGAGA"
$ws.Cells.Item(38, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCAGCAAGAGATCGTCGGCAGCGTC"
$ws.Cells.Item(39, 2).Value = "F1094-CTGTGTCGAT"
$ws.Cells.Item(39, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTGTGTCGATTCGTCGGCAGCGTC"
$ws.Cells.Item(40, 2).Value = "F1095-GACGACTCAA"
$ws.Cells.Item(40, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGACGACTCAATCGTCGGCAGCGTC"
$ws.Cells.Item(41, 2).Value = "F1096-CCTTCAAGTC"
$ws.Cells.Item(41, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCCTTCAAGTCTCGTCGGCAGCGTC"
$ws.Cells.Item(42, 2).Value = "F1097-AAGGTAGTTG"
$ws.Cells.Item(42, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAAGGTAGTTGTCGTCGGCAGCGTC"
$ws.Cells.Item(43, 2).Value = "F1098-TGACAGTGAG"
$ws.Cells.Item(43, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGACAGTGAGTCGTCGGCAGCGTC"
$ws.Cells.Item(44, 2).Value = "F1099-ATGGTGTGTT"
$ws.Cells.Item(44, 3).Value = "AATGATACGGCGACCACCGAGATCTACACATGGTGTGTTTCGTCGGCAGCGTC"
$ws.Cells.Item(45, 2).Value = "F1100-GTACGAACAA"
$ws.Cells.Item(45, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTACGAACAATCGTCGGCAGCGTC"
$ws.Cells.Item(46, 2).Value = "F1101-CATCACGTAG"
$ws.Cells.Item(46, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCATCACGTAGTCGTCGGCAGCGTC"
$ws.Cells.Item(47, 2).Value = "F1102-ATCCTTCGAC"
$ws.Cells.Item(47, 3).Value = "AATGATACGGCGACCACCGAGATCTACACATCCTTCGACTCGTCGGCAGCGTC"
$ws.Cells.Item(48, 2).Value = "F1103-CTGTAGCAGA"
$ws.Cells.Item(48, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTGTAGCAGATCGTCGGCAGCGTC"
$ws.Cells.Item(49, 2).Value = "F1104-CTCAACAGTG"
$ws.Cells.Item(49, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTCAACAGTGTCGTCGGCAGCGTC"
$ws.Cells.Item(50, 2).Value = "F1105-CTTCAACACC"
$ws.Cells.Item(50, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTTCAACACCTCGTCGGCAGCGTC"
$ws.Cells.Item(51, 2).Value = "F1106-CAGACACGTT"
$ws.Cells.Item(51, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCAGACACGTTTCGTCGGCAGCGTC"
$ws.Cells.Item(52, 2).Value = "F1107-ATGCAAGACC"
$ws.Cells.Item(52, 3).Value = "AATGATACGGCGACCACCGAGATCTACACATGCAAGACCTCGTCGGCAGCGTC"
$ws.Cells.Item(53, 2).Value = "F1108-CTGTGAAGGA"
$ws.Cells.Item(53, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTGTGAAGGATCGTCGGCAGCGTC"
$ws.Cells.Item(54, 2).Value = "F1109-ACATCACTGG"
$ws.Cells.Item(54, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACATCACTGGTCGTCGGCAGCGTC"
$ws.Cells.Item(55, 2).Value = "F1110-TACGACGTAG"
$ws.Cells.Item(55, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTACGACGTAGTCGTCGGCAGCGTC"
$ws.Cells.Item(56, 2).Value = "F1111-AACAGACTGG"
$ws.Cells.Item(56, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAACAGACTGGTCGTCGGCAGCGTC"
$ws.Cells.Item(57, 2).Value = "F1112-CCACATCACT"
$ws.Cells.Item(57, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCCACATCACTTCGTCGGCAGCGTC"
$ws.Cells.Item(58, 2).Value = "F1113-AAGTCACATC"
$ws.Cells.Item(58, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAAGTCACATCTCGTCGGCAGCGTC"
$ws.Cells.Item(59, 2).Value = "F1114-GCAAGATCTC"
$ws.Cells.Item(59, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGCAAGATCTCTCGTCGGCAGCGTC"
$ws.Cells.Item(60, 2).Value = "F1115-TGTCGAACAC"
$ws.Cells.Item(60, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGTCGAACACTCGTCGGCAGCGTC"
$ws.Cells.Item(61, 2).Value = "F1116-TCTACAGTCG"
$ws.Cells.Item(61, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCTACAGTCGTCGTCGGCAGCGTC"
$ws.Cells.Item(62, 2).Value = "F1117-TCTTCAGACT"
$ws.Cells.Item(62, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCTTCAGACTTCGTCGGCAGCGTC"
$ws.Cells.Item(63, 2).Value = "F1118-GTCCTTCAAG"
$ws.Cells.Item(63, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTCCTTCAAGTCGTCGGCAGCGTC"
$ws.Cells.Item(64, 2).Value = "F1119-CCTAGCAGTA"
$ws.Cells.Item(64, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCCTAGCAGTATCGTCGGCAGCGTC"
$ws.Cells.Item(65, 2).Value = "F1120-TTGCAGCAAG"
$ws.Cells.Item(65, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTTGCAGCAAGTCGTCGGCAGCGTC"
$ws.Cells.Item(66, 2).Value = "F1121-AAGAAGCTGA"
$ws.Cells.Item(66, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAAGAAGCTGATCGTCGGCAGCGTC"
$ws.Cells.Item(67, 2).Value = "F1122-CATGCAGCAT"
$ws.Cells.Item(67, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCATGCAGCATTCGTCGGCAGCGTC"
$ws.Cells.Item(68, 2).Value = "F1123-CAGAGTTGGA"
$ws.Cells.Item(68, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCAGAGTTGGATCGTCGGCAGCGTC"
$ws.Cells.Item(69, 2).Value = "F1124-AGAAGAGGTC"
$ws.Cells.Item(69, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAGAAGAGGTCTCGTCGGCAGCGTC"
$ws.Cells.Item(70, 2).Value = "F1125-ACTTCCAACC"
$ws.Cells.Item(70, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACTTCCAACCTCGTCGGCAGCGTC"
$ws.Cells.Item(71, 2).Value = "F1126-GGACTACGTT"
$ws.Cells.Item(71, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGGACTACGTTTCGTCGGCAGCGTC"
$ws.Cells.Item(72, 2).Value = "F1127-TGACTCCTTC"
$ws.Cells.Item(72, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGACTCCTTCTCGTCGGCAGCGTC"
$ws.Cells.Item(73, 2).Value = "F1128-GAGGTACAGT"
$ws.Cells.Item(73, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGAGGTACAGTTCGTCGGCAGCGTC"
$ws.Cells.Item(74, 2).Value = "F1129-CAAGTAGATG"
$ws.Cells.Item(74, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCAAGTAGATGTCGTCGGCAGCGTC"
$ws.Cells.Item(75, 2).Value = "F1130-CTTCTCTGTT"
$ws.Cells.Item(75, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTTCTCTGTTTCGTCGGCAGCGTC"
$ws.Cells.Item(76, 2).Value = "F1131-TAGACCAGGT"
$ws.Cells.Item(76, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTAGACCAGGTTCGTCGGCAGCGTC"
$ws.Cells.Item(77, 2).Value = "F1132-TACAAGAGGT"
$ws.Cells.Item(77, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTACAAGAGGTTCGTCGGCAGCGTC"
$ws.Cells.Item(78, 2).Value = "F1133-CTTCACACCT"
$ws.Cells.Item(78, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCTTCACACCTTCGTCGGCAGCGTC"
$ws.Cells.Item(79, 2).Value = "F1134-TTCCACTGTG"
$ws.Cells.Item(79, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTTCCACTGTGTCGTCGGCAGCGTC"
$ws.Cells.Item(80, 2).Value = "F1135-TACTTGCACC"
$ws.Cells.Item(80, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTACTTGCACCTCGTCGGCAGCGTC"
$ws.Cells.Item(81, 2).Value = "F1136-TCTCATGGAT"
$ws.Cells.Item(81, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCTCATGGATTCGTCGGCAGCGTC"
$ws.Cells.Item(82, 2).Value = "F1137-ACGATGGTAC"
$ws.Cells.Item(82, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACGATGGTACTCGTCGGCAGCGTC"
$ws.Cells.Item(83, 2).Value = "F1138-TCAGAACACT"
$ws.Cells.Item(83, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCAGAACACTTCGTCGGCAGCGTC"
$ws.Cells.Item(84, 2).Value = "F1139-TGGTCCTTGA"
$ws.Cells.Item(84, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGGTCCTTGATCGTCGGCAGCGTC"
$ws.Cells.Item(85, 2).Value = "F1140-AAGTCTCCAA"
$ws.Cells.Item(85, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAAGTCTCCAATCGTCGGCAGCGTC"
$ws.Cells.Item(86, 2).Value = "F1141-TGAGACGTTG"
$ws.Cells.Item(86, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTGAGACGTTGTCGTCGGCAGCGTC"
$ws.Cells.Item(87, 2).Value = "F1142-TAGACCTTCT"
$ws.Cells.Item(87, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTAGACCTTCTTCGTCGGCAGCGTC"
$ws.Cells.Item(88, 2).Value = "F1143-ACGTAGACTA"
$ws.Cells.Item(88, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACGTAGACTATCGTCGGCAGCGTC"
$ws.Cells.Item(89, 2).Value = "F1144-TCTCTCTCTT"
$ws.Cells.Item(89, 3).Value = "AATGATACGGCGACCACCGAGATCTACACTCTCTCTCTTTCGTCGGCAGCGTC"
$ws.Cells.Item(90, 2).Value = "F1145-AACACTGGAA"
$ws.Cells.Item(90, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAACACTGGAATCGTCGGCAGCGTC"
$ws.Cells.Item(91, 2).Value = "F1146-AACCAGAGCT"
$ws.Cells.Item(91, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAACCAGAGCTTCGTCGGCAGCGTC"
$ws.Cells.Item(92, 2).Value = "F1147-GGATGTCGAT"
$ws.Cells.Item(92, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGGATGTCGATTCGTCGGCAGCGTC"
$ws.Cells.Item(93, 2).Value = "F1148-AACCTTGTGA"
$ws.Cells.Item(93, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAACCTTGTGATCGTCGGCAGCGTC"
$ws.Cells.Item(94, 2).Value = "F1149-AAGGATGCAA"
$ws.Cells.Item(94, 3).Value = "AATGATACGGCGACCACCGAGATCTACACAAGGATGCAATCGTCGGCAGCGTC"
$ws.Cells.Item(95, 2).Value = "F1150-GTAGCTAGTG"
$ws.Cells.Item(95, 3).Value = "AATGATACGGCGACCACCGAGATCTACACGTAGCTAGTGTCGTCGGCAGCGTC"
$ws.Cells.Item(96, 2).Value = "F1151-ACAGTTCCTA"
$ws.Cells.Item(96, 3).Value = "AATGATACGGCGACCACCGAGATCTACACACAGTTCCTATCGTCGGCAGCGTC"
$ws.Cells.Item(97, 2).Value = "F1152-CATGTACGTC"
$ws.Cells.Item(97, 3).Value = "AATGATACGGCGACCACCGAGATCTACACCATGTACGTCTCGTCGGCAGCGTC"
